$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 00:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1815102
$ws.Range("C4").Value = 21572
$ws.Range("D4").Value = 528370
$ws.Range("E4").Value = 1181211
$ws.Range("G4").Value = 979
$ws.Range("H4").Value = 105521

# Brasil (row 5)
$ws.Range("D5").Value = 205371
$ws.Range("E5").Value = 236124

# Alemania (row 11)
$ws.Range("B11").Value = 183294
$ws.Range("C11").Value = 275
$ws.Range("E11").Value = 9794

# Canada (row 17)
$ws.Range("B17").Value = 90166
$ws.Range("C17").Value = 748
$ws.Range("D17").Value = 48065
$ws.Range("E17").Value = 35028

# Republica Dominicana moves above Japon & Austria (new data lands on row 44,
# Japon and Austria's previous data shift down to rows 45 and 46)
$ws.Range("A44").Value = "Republica Dominicana"
$ws.Range("A45").Value = "Japon"
$ws.Range("A46").Value = "Austria"
$ws.Range("B44").Value = 16908
$ws.Range("C44").Value = 377
$ws.Range("D44").Value = 9557
$ws.Range("E44").Value = 6853
$ws.Range("G44").Value = 10
$ws.Range("H44").Value = 498

$ws.Range("B45").Value = 16719
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 14254
$ws.Range("E45").Value = 1591
$ws.Range("H45").Value = 874

$ws.Range("B46").Value = 16685
$ws.Range("C46").Value = 30
$ws.Range("D46").Value = 15520
$ws.Range("E46").Value = 497
$ws.Range("H46").Value = 668

# Guinea (row 77)
$ws.Range("B77").Value = 3706
$ws.Range("C77").Value = 50
$ws.Range("D77").Value = 2030
$ws.Range("E77").Value = 1653
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 23

# Maldivas (row 98)
$ws.Range("D98").Value = 406
$ws.Range("E98").Value = 1261

# Uganda moves above Ruanda, Isla de Man, Mauricio (new data lands on row 146,
# those three countries' previous data shift down to rows 147-149)
$ws.Range("A146").Value = "Uganda"
$ws.Range("A147").Value = "Ruanda"
$ws.Range("A148").Value = "Isla de Man"
$ws.Range("A149").Value = "Mauricio"
$ws.Range("B146").Value = 413
$ws.Range("C146").Value = 84
$ws.Range("D146").Value = 72
$ws.Range("E146").Value = 341
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

$ws.Range("B147").Value = 359
$ws.Range("C147").Value = 4
$ws.Range("D147").Value = 250
$ws.Range("E147").Value = 108
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 1

$ws.Range("B148").Value = 336
$ws.Range("D148").Value = 309
$ws.Range("H148").Value = 24

$ws.Range("B149").Value = 335
$ws.Range("D149").Value = 322
$ws.Range("E149").Value = 3
$ws.Range("H149").Value = 10
